$d = $word.ActiveDocument

# HPCDATAMGM-2107: merge Fix/ed/ the issue... runs into one (same text)
$d.Content.Find.Execute("Fixed the issue of the Permissions Summary page of the DME web application showing only the permissions for the first user when permissions for a set of files selected from the search results screen are applied to a list of users. ", $true, $false, $false, $false, $false, $true, 1, $false, "Fixed the issue of the Permissions Summary page of the DME web application showing only the permissions for the first user when permissions for a set of files selected from the search results screen are applied to a list of users. ", 2) | Out-Null

# HPCDATAMGM-2117, 2089: Add -> Added, data object -> data objects
$d.Content.Find.Execute("Add the ability for DME to locate data object residing in a ‘linked’ archive  ", $true, $false, $false, $false, $false, $true, 1, $false, "Added the ability for DME to locate data objects residing in a ‘linked’ archive  ", 2) | Out-Null

# HPCDATAMGM-2110: merge Create/d/ DB indexes... runs into one (same text)
$d.Content.Find.Execute("Created DB indexes on the DME tables recommended by the Oracle Enterprise Manager to improve performance. ", $true, $false, $false, $false, $false, $true, 1, $false, "Created DB indexes on the DME tables recommended by the Oracle Enterprise Manager to improve performance. ", 2) | Out-Null

# HPCDATAMGM-1910: merge ": Remediated ..." runs into one (same text)
$d.Content.Find.Execute(": Remediated the security vulnerability flagged by the Nessus scan on the libcrypto.so library copied by the Aspera Connect installation on the DME server.", $true, $false, $false, $false, $false, $true, 1, $false, ": Remediated the security vulnerability flagged by the Nessus scan on the libcrypto.so library copied by the Aspera Connect installation on the DME server.", 2) | Out-Null

# HPCDATAMGM-2100...: merge ": Upgrade/d ..." runs into one (same text)
$d.Content.Find.Execute(": Upgraded open-source libraries used in DME that have been identified as having security vulnerabilities by the Nessus scan.", $true, $false, $false, $false, $false, $true, 1, $false, ": Upgraded open-source libraries used in DME that have been identified as having security vulnerabilities by the Nessus scan.", 2) | Out-Null
